$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.619.14'
$ws.Cells.Item(2, 5).Value = '  +3.30%  '

$ws.Cells.Item(3, 4).Value = '2.446.01'
$ws.Cells.Item(3, 5).Value = '  +2.04%  '

$ws.Cells.Item(4, 5).Value = '  -0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '578.72'
$ws.Cells.Item(5, 5).Value = '  +2.82%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '145.53'
$ws.Cells.Item(6, 5).Value = '  +3.36%  '

$ws.Cells.Item(7, 5).Value = '  +0.03%  '

$ws.Cells.Item(8, 5).Value = '  +0.60%  '

$ws.Cells.Item(9, 4).Value = '2.444.38'
$ws.Cells.Item(9, 5).Value = '  +1.76%  '

$ws.Cells.Item(10, 5).Value = '  +2.55%  '

$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.160'
$ws.Cells.Item(11, 5).Value = '  +1.09%  '

$ws.Cells.Item(12, 2).Value = 'Toncoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.23'
$ws.Cells.Item(12, 5).Value = '  +1.41%  '

$ws.Cells.Item(13, 2).Value = 'Cardano'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.353'
$ws.Cells.Item(13, 5).Value = '  +3.78%  '

$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '28.44'
$ws.Cells.Item(14, 5).Value = '  +9.15%  '

$ws.Cells.Item(15, 2).Value = 'ShibaInu'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000178'
$ws.Cells.Item(15, 5).Value = '  +6.06%  '

$ws.Cells.Item(16, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(16, 4).Value = '2.889.92'
$ws.Cells.Item(16, 5).Value = '  +1.89%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '62.561.38'
$ws.Cells.Item(17, 5).Value = '  +3.68%  '

$ws.Cells.Item(18, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(18, 4).Value = '0.0₅0111'
$ws.Cells.Item(18, 5).Value = '  +295.77%  '

$ws.Cells.Item(19, 4).Value = '2.442.43'
$ws.Cells.Item(19, 5).Value = '  +1.51%  '

$ws.Cells.Item(20, 5).Value = '  -1.45%  '

$ws.Cells.Item(21, 5).Value = '  +2.60%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '325.94'
$ws.Cells.Item(22, 5).Value = '  +0.78%  '

$ws.Cells.Item(23, 5).Value = '  +1.11%  '

$ws.Cells.Item(24, 5).Value = '  +11.16%  '

$ws.Cells.Item(25, 5).Value = '  +0.01%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '65.31'
$ws.Cells.Item(26, 5).Value = '  +0.50%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '645.88'
$ws.Cells.Item(27, 5).Value = '  +14.88%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.15'
$ws.Cells.Item(28, 5).Value = '  +14.99%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.56'
$ws.Cells.Item(29, 5).Value = '  +6.38%  '

$ws.Cells.Item(30, 5).Value = '  +4.82%  '

$ws.Cells.Item(31, 4).Value = '2.558.81'
$ws.Cells.Item(31, 5).Value = '  +1.83%  '

$ws.Cells.Item(32, 5).Value = '  +1.30%  '

$ws.Cells.Item(33, 5).Value = '  +6.82%  '

$ws.Cells.Item(34, 5).Value = '  +3.44%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.138'
$ws.Cells.Item(35, 5).Value = '  +6.11%  '

$ws.Cells.Item(36, 5).Value = '  +2.19%  '

$ws.Cells.Item(37, 5).Value = '  +0.09%  '

$ws.Cells.Item(38, 5).Value = '  +3.44%  '

$ws.Cells.Item(39, 2).Value = 'Monero'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '153.94'
$ws.Cells.Item(39, 5).Value = '  +1.14%  '

$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.47'
$ws.Cells.Item(40, 5).Value = '  +6.78%  '

$ws.Cells.Item(41, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.373'
$ws.Cells.Item(41, 5).Value = '  +0.79%  '

$ws.Cells.Item(42, 5).Value = '  +1.84%  '

$ws.Cells.Item(43, 5).Value = '  +8.87%  '

$ws.Cells.Item(44, 5).Value = '  +5.09%  '

$ws.Cells.Item(45, 5).Value = '  +2.06%  '

$ws.Cells.Item(47, 5).Value = '  +28.06%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '144.34'
$ws.Cells.Item(48, 5).Value = '  +2.19%  '

$ws.Cells.Item(49, 5).Value = '  +1.58%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '20.56'
$ws.Cells.Item(50, 5).Value = '  +7.05%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.603'
$ws.Cells.Item(51, 5).Value = '  +2.31%  '

